# [IMP] Odoo 14 bank statement import
# Fill in the (previously blank) Currency column for the sample data rows,
# and nudge the cosmetic view state to match what Excel would save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Currency column (E) values for the two data rows.
$ws.Range("E2").Value = "KWD"
$ws.Range("E3").Value = "KWD"

# Match the column-3 ("Partner") width tweak from the diff.
$ws.Columns.Item(3).ColumnWidth = 23.23

# Match the final cell selection recorded in the sheet view.
$ws.Range("E5").Select()
